# Split the single "SectionHeadnote" paragraph that holds the whole
# section/case blob of text into the proper sequence of resource/case/
# section paragraphs (each with its own style), plus the table-of-contents
# bookmarks the document's TOC hyperlinks already point at.

$d = $word.ActiveDocument

# Locate the paragraph that currently holds the big concatenated blob of
# text ("1.1Case of the District Number 1This is the body of case 1. ...")
# by searching for a short, unique anchor at its very start.
$hit = $d.Content
$found = $hit.Find.Execute("1.1Case of the District Number 1", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target SectionHeadnote paragraph"
}

# Grow the found hit out to the whole paragraph so we replace all of its
# text (the paragraph mark itself is left alone).
$target = $d.Range($hit.Start, $hit.Start)
$target.Expand(4) | Out-Null

# New text for each resulting paragraph, in document order.
$parts = @(
    "What is a corporation?",
    "1.1",
    "Case of the District Number 1",
    "",
    "This is the body of case 1.",
    "1.2",
    "Case of the District Number 2",
    "This is an annotatable resource in the casebook.`n",
    "highlighted: content to highlight; elided: content to elide; replaced: content to replace; commented: content to comment; highlighted2: second highlight content;`n",
    "2",
    "Section Two",
    "This is the second chapter of the casebook.`n"
)

# Paragraph style to apply to each of the new paragraphs, same order.
$styles = @(
    "SectionHeadnote",
    "ResourceNumber",
    "ResourceTitle",
    "ResourceHeadnote",
    "CaseText",
    "ResourceNumber",
    "ResourceTitle",
    "ResourceHeadnote",
    "CaseText",
    "SectionNumber",
    "SectionTitle",
    "SectionHeadnote"
)

# Replacing the text with a single string containing carriage returns
# splits it into that many paragraphs in one shot (each new paragraph
# inherits the original SectionHeadnote style to start with).
$target.Text = [string]::Join("`r", $parts)

# The new paragraphs are the last $parts.Length paragraphs in the document.
$totalParas = $d.Paragraphs.Count
$startIndex = $totalParas - $parts.Length + 1

for ($i = 0; $i -lt $parts.Length; $i++) {
    $d.Paragraphs.Item($startIndex + $i).Style = $styles[$i]
}

# Re-create the TOC bookmarks (_auto_toc_2/_auto_toc_3/_auto_toc_4) around
# the resource/section numbers, matching what the existing TOC hyperlinks
# already point at. Re-fetch each target range by absolute position so the
# Bookmarks.Add call sees a freshly bound Range.
function Add-NumberBookmark($partIndex, $name) {
    $p = $d.Paragraphs.Item($startIndex + $partIndex)
    $s = $p.Range.Start
    $e = $p.Range.End
    $bmRange = $d.Range($s, $e)
    $d.Bookmarks.Add($name, $bmRange) | Out-Null
}

Add-NumberBookmark 1 "_auto_toc_2"   # "1.1"
Add-NumberBookmark 5 "_auto_toc_3"   # "1.2"
Add-NumberBookmark 9 "_auto_toc_4"   # "2"
